$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105, shifting existing rows 105-212 down to 106-213
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row 105 with the new data record
$ws.Range("A105").Value = 6
$ws.Range("B105").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C105").Value = "Metropolitana"
$ws.Range("D105").Value = 44705
$ws.Range("E105").Value = 13
$ws.Range("F105").Value = 100112001
$ws.Range("G105").Value = "Berenjena"
$ws.Range("H105").Value = "Sin especificar"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 200
$ws.Range("K105").Value = 4000
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = 4400
$ws.Range("N105").Value = "`$/caja 50 unidades"
$ws.Range("O105").Value = "Región de Arica y Parinacota"
$ws.Range("P105").Value = 88
$ws.Range("Q105").Value = 50
$ws.Range("R105").Value = "Hortaliza"

# Ensure the date style (numFmt for dates) used by column D is applied to the new D105 cell
$ws.Range("D105").NumberFormat = $ws.Range("D106").NumberFormat
